$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# 1. Insert a new column before column N; this shifts N..AC one column to the
#    right (N->O, O->P, ... AC->AD), carrying along all the existing
#    command-list data for macro/mail/number/pdf/rdbms/redis/sms/sound/ssh/
#    step/web/webalert/webcookie/ws/ws.async/xml.
$ws.Columns("N:N").Insert()

# 2. Populate the freshly inserted column N with the new "localdb" target
#    and its six commands.
$ws.Range("N1").Value = "localdb"
$ws.Range("N2").Value = "cloneTable(var,source,target)"
$ws.Range("N3").Value = "dropTables(var,tables)"
$ws.Range("N4").Value = "exportCSV(sql,output)"
$ws.Range("N5").Value = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value = "purge(var)"
$ws.Range("N7").Value = "runSQLs(var,sqls)"

# 3. The "target" list in column A (rows 2-29) is alphabetically sorted;
#    insert "localdb" between "json" and "macro" by re-writing the list,
#    which pushes every following entry down one row (list now spans
#    rows 2-30).
$ws.Range("A14").Value = "localdb"
$ws.Range("A15").Value = "macro"
$ws.Range("A16").Value = "mail"
$ws.Range("A17").Value = "number"
$ws.Range("A18").Value = "pdf"
$ws.Range("A19").Value = "rdbms"
$ws.Range("A20").Value = "redis"
$ws.Range("A21").Value = "sms"
$ws.Range("A22").Value = "sound"
$ws.Range("A23").Value = "ssh"
$ws.Range("A24").Value = "step"
$ws.Range("A25").Value = "web"
$ws.Range("A26").Value = "webalert"
$ws.Range("A27").Value = "webcookie"
$ws.Range("A28").Value = "ws"
$ws.Range("A29").Value = "ws.async"
$ws.Range("A30").Value = "xml"

# 4. Update the defined names so the existing ranges point at their new
#    (shifted) columns, extend "target" to include the new row, and add
#    the brand-new "localdb" named range.
$wb.Names.Item("macro").RefersTo = "='#system'!`$O`$2:`$O`$4"
$wb.Names.Item("mail").RefersTo = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("sms").RefersTo = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ssh").RefersTo = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$21"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
